$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the asset code, name and location for the single data row.
# "02.02.0004" looks like a date to Excel's auto-detection, so enter it
# as a text formula first and then paste-special as a value; this keeps
# the result as a plain text/shared-string cell without touching the
# cell's number format / style (which a direct NumberFormat="@" change
# would otherwise introduce).
$ws.Range("B2").Formula = '="02.02.0004"'
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)

$ws.Range("C2").Value = "Kantor Prabubima Tech"
$ws.Range("I2").Value = "Pabuaran"

# Update the view: scroll back to column A and move the selection to I3
# (matches the saved sheetView: no topLeftCell scroll, selection at I3)
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("I3").Select()
